$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: new client record (Process = "no" -> will be skipped on load) ---
# Fill the descriptive columns first, then the Process flag last (matches
# natural data-entry order: you decide Process after reviewing the row).
$ws.Range("B9").Value = "Center for Service of National Police Units of Ukraine"
$ws.Range("C9").Value = "https://mvs.gov.ua/en/contacts/national-police-ukraine`n"
$ws.Range("D9").Value = "anketa.bezvisty@mvs.gov.ua"
$ws.Range("E9").Value = "Anketa Bezvisty"
$ws.Range("A9").Value = "no"

# Hyperlinks for the new row's website + email cells
$ws.Hyperlinks.Add($ws.Range("C9"), "https://mvs.gov.ua/en/contacts/national-police-ukraine")
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:anketa.bezvisty@mvs.gov.ua")

# Website column keeps the hyperlink look; recipient_email column does not
$ws.Range("C9").Style = "Hyperlink"
$ws.Range("D9").Style = "Normal"

# Undo the engine's auto row-height bump triggered by the newline in C9
$ws.Rows.Item(9).AutoFit()

# Row 8's recipient_email cell (D8) loses its hyperlink styling (still a live link)
$ws.Range("D8").Style = "Normal"

# Leftover formatted-but-empty row below the new data (matches manual entry artifact)
$ws.Range("C10").Style = "Hyperlink"
$ws.Range("D10").Style = "Hyperlink"

# --- Column widths ---
$ws.Columns.Item(3).ColumnWidth = 31.25
$ws.Columns.Item(4).ColumnWidth = 26.25

# --- Selection moved to C14 ---
[void]$ws.Range("C14").Select()

Write-Output "edit complete"
